$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header for new "Speed" column (F34)
$ws.Range("F34").Value = "Speed"

# Update Health (column D) values for rows 35-45 and add Speed (column F) values
$ws.Range("D35").Value = 50
$ws.Range("F35").Value = 1

$ws.Range("D36").Value = 100
$ws.Range("F36").Value = 1

$ws.Range("D37").Value = 1000
$ws.Range("F37").Value = 1

$ws.Range("D38").Value = 100
$ws.Range("F38").Value = 2

$ws.Range("D39").Value = 80
$ws.Range("F39").Value = 1

$ws.Range("D40").Value = 60
$ws.Range("F40").Value = 1

$ws.Range("D41").Value = 50
$ws.Range("F41").Value = 3

$ws.Range("D42").Value = 50
$ws.Range("F42").Value = 2

$ws.Range("D43").Value = 150
$ws.Range("F43").Value = 1

$ws.Range("D44").Value = 80
$ws.Range("F44").Value = 1

$ws.Range("D45").Value = 40
$ws.Range("F45").Value = 1

# Update the view: scroll position and active selection
$ws.Activate()
$ws.Range("F46").Select()
